# Updated cryptos list (price + volume refresh, and a row swap)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Temporarily force Text format on the Price/Volume columns so that
# numeric-looking strings (e.g. '243.97', '30.336.02') are stored as
# literal text, matching the source data feed's formatting.
$priceVolumeRange = $ws.Range('D2:E51')
$priceVolumeRange.NumberFormat = '@'

$ws.Range('D2').Value = '30.336.02'
$ws.Range('E2').Value = '  +0.33%  '
$ws.Range('D3').Value = '1.933.60'
$ws.Range('E3').Value = '  +0.25%  '
$ws.Range('E4').Value = '  -0.09%  '
$ws.Range('D5').Value = '0.7535'
$ws.Range('E5').Value = '  +5.88%  '
$ws.Range('D6').Value = '243.97'
$ws.Range('E6').Value = '  -1.90%  '
$ws.Range('D7').Value = '1.001'
$ws.Range('E7').Value = '  -0.08%  '
$ws.Range('D8').Value = '27.97'
$ws.Range('E8').Value = '  +2.56%  '
$ws.Range('D9').Value = '0.3185'
$ws.Range('E9').Value = '  -0.81%  '
$ws.Range('D10').Value = '0.07029'
$ws.Range('E10').Value = '  -0.94%  '
$ws.Range('D11').Value = '0.7809'
$ws.Range('E11').Value = '  -1.37%  '
$ws.Range('D12').Value = '0.08029'
$ws.Range('E12').Value = '  +0.06%  '
$ws.Range('D13').Value = '1.924.14'
$ws.Range('E13').Value = '  -0.24%  '
$ws.Range('D14').Value = '5.400'
$ws.Range('E14').Value = '  +0.62%  '
$ws.Range('D15').Value = '93.38'
$ws.Range('E15').Value = '  -1.49%  '
$ws.Range('E16').Value = '  -1.61%  '
$ws.Range('D17').Value = '30.338.26'
$ws.Range('E17').Value = '  +0.28%  '
$ws.Range('B18').Value = 'Uniswap'
$ws.Range('C18').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D18').Value = '6.058'
$ws.Range('E18').Value = '  +5.38%  '
$ws.Range('B19').Value = 'BitcoinCash'
$ws.Range('C19').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D19').Value = '252.77'
$ws.Range('E19').Value = '  -0.91%  '
$ws.Range('D20').Value = '0.000007967'
$ws.Range('E20').Value = '  -0.73%  '
$ws.Range('D21').Value = '2.180.53'
$ws.Range('E21').Value = '  -0.05%  '
$ws.Range('D22').Value = '0.9994'
$ws.Range('E22').Value = '  -0.24%  '
$ws.Range('D23').Value = '1.002'
$ws.Range('E23').Value = '  -0.06%  '
$ws.Range('D24').Value = '6.707'
$ws.Range('E24').Value = '  -1.62%  '
$ws.Range('D25').Value = '9.521'
$ws.Range('E25').Value = '  -0.28%  '
$ws.Range('D26').Value = '164.48'
$ws.Range('E26').Value = '  -1.02%  '
$ws.Range('D27').Value = '19.09'
$ws.Range('E27').Value = '  +0.22%  '
$ws.Range('E28').Value = '  +3.07%  '
$ws.Range('D29').Value = '2.223'
$ws.Range('E29').Value = '  -2.19%  '
$ws.Range('D30').Value = '1.367'
$ws.Range('E30').Value = '  +0.72%  '
$ws.Range('D31').Value = '1.527'
$ws.Range('E31').Value = '  -0.12%  '
$ws.Range('D32').Value = '4.409'
$ws.Range('E32').Value = '  +0.42%  '
$ws.Range('E33').Value = '  +0.04%  '
$ws.Range('D34').Value = '0.05228'
$ws.Range('E34').Value = '  +1.19%  '
$ws.Range('D35').Value = '1.323'
$ws.Range('E35').Value = '  +4.60%  '
$ws.Range('D36').Value = '0.7546'
$ws.Range('D37').Value = '2.790'
$ws.Range('E37').Value = '  +0.77%  '
$ws.Range('D38').Value = '0.01950'
$ws.Range('E38').Value = '  -0.16%  '
$ws.Range('D39').Value = '2.807'
$ws.Range('E39').Value = '  +0.14%  '
$ws.Range('D40').Value = '6.617'
$ws.Range('E40').Value = '  +4.32%  '
$ws.Range('D41').Value = '78.82'
$ws.Range('E41').Value = '  +1.48%  '
$ws.Range('D42').Value = '0.4489'
$ws.Range('E42').Value = '  +0.21%  '
$ws.Range('D43').Value = '1.974'
$ws.Range('E43').Value = '  +0.05%  '
$ws.Range('D44').Value = '1.001'
$ws.Range('E44').Value = '  -0.03%  '
$ws.Range('D45').Value = '0.8358'
$ws.Range('E45').Value = '  -1.19%  '
$ws.Range('D46').Value = '9.943'
$ws.Range('E46').Value = '  +2.79%  '
$ws.Range('E47').Value = '  +0.98%  '
$ws.Range('D48').Value = '7.605'
$ws.Range('E48').Value = '  +2.41%  '
$ws.Range('D49').Value = '37.90'
$ws.Range('E49').Value = '  +4.25%  '
$ws.Range('D50').Value = '983.25'
$ws.Range('E50').Value = '  +7.53%  '
$ws.Range('D51').Value = '0.1224'
$ws.Range('E51').Value = '  +7.98%  '

# Restore the original (default) cell style now that the text values are set.
$priceVolumeRange.Style = 'Normal'
